$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Task1"
$ws.Range("C5").Clear()
$ws.Range("P7").Value = "Task"
$ws.Range("P8").Formula = '= N8+COUNTIF($L$8:L8,L8)-1'
$ws.Range("P9").Formula = '= N9+COUNTIF($L$8:L9,L9)-1'
$ws.Range("P10").Formula = '= N10+COUNTIF($L$8:L10,L10)-1'
$ws.Range("P11").Formula = '= N11+COUNTIF($L$8:L11,L11)-1'
$ws.Range("P12").Formula = '= N12+COUNTIF($L$8:L12,L12)-1'
$ws.Range("P13").Formula = '= N13+COUNTIF($L$8:L13,L13)-1'
$ws.Range("P14").Formula = '= N14+COUNTIF($L$8:L14,L14)-1'
